$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2136
$ws.Range("J17").Value = 2136
$ws.Range("L17").Value = 6408
$ws.Range("N17").Value = -6744
# Row 19
$ws.Range("H19").Value = 1296.8334
$ws.Range("J19").Value = 1342
$ws.Range("L19").Value = 1342
$ws.Range("N19").Value = -1692
# Row 112
$ws.Range("H112").Value = 4611.1113
$ws.Range("I112").Value = 1000
$ws.Range("K112").Value = 3000
$ws.Range("M112").Value = -1892
# Row 125
$ws.Range("H125").Value = 1330
$ws.Range("J125").Value = 1350
$ws.Range("L125").Value = 12150
$ws.Range("N125").Value = -17070

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 12
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
# Row 32
$ws.Range("H32").Value = 2548.3157
$ws.Range("I32").Value = 2548.3157
$ws.Range("K32").Value = 2548.3157
$ws.Range("M32").Value = -2261.3157
# Row 76
$ws.Range("H76").Value = 46661.668
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 46661.668
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 46661.668
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -47337.668
# Row 79
$ws.Range("H79").Value = 46661.668
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 46661.668
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 46661.668
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -49001.668
# Row 109
$ws.Range("H109").Value = 28969
$ws.Range("J109").Value = 28969
$ws.Range("L109").Value = 28969
$ws.Range("N109").Value = -31743
# Row 132
$ws.Range("H132").Value = 4921.8335
$ws.Range("I132").Value = 4911.4287
$ws.Range("J132").Value = 4994.6665
$ws.Range("K132").Value = 14734.2861
$ws.Range("L132").Value = 14983.9995
$ws.Range("M132").Value = -12204.2861
$ws.Range("N132").Value = -20043.9995

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2931.5
$ws.Range("J20").Value = 2499
$ws.Range("L20").Value = 2499
$ws.Range("N20").Value = -2993
# Row 22
$ws.Range("H22").Value = 368
$ws.Range("I22").Value = 458.25
$ws.Range("K22").Value = 458.25
$ws.Range("M22").Value = -285.25
# Row 86
$ws.Range("H86").Value = 6235.6895
$ws.Range("I86").Value = 6797.522
$ws.Range("J86").Value = 4082
$ws.Range("K86").Value = 6797.522
$ws.Range("L86").Value = 4082
$ws.Range("M86").Value = -5674.522
$ws.Range("N86").Value = -6328
# Row 89
$ws.Range("H89").Value = 6235.6895
$ws.Range("I89").Value = 6797.522
$ws.Range("J89").Value = 4082
$ws.Range("K89").Value = 33987.61
$ws.Range("L89").Value = 20410
$ws.Range("M89").Value = -28371.61
$ws.Range("N89").Value = -31642

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 25
$ws.Range("H25").Value = 24999
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
# Row 31
$ws.Range("H31").Value = 2072.0588
$ws.Range("I31").Value = 2101.8667
$ws.Range("K31").Value = 2101.8667
$ws.Range("M31").Value = -1806.8667
# Row 34
$ws.Range("H34").Value = 2072.0588
$ws.Range("I34").Value = 2101.8667
$ws.Range("K34").Value = 2101.8667
$ws.Range("M34").Value = -1899.8667
# Row 99
$ws.Range("H99").Value = 5507.5
$ws.Range("I99").Value = 3998
$ws.Range("K99").Value = 3998
$ws.Range("M99").Value = -2500
# Row 126
$ws.Range("H126").Value = 5507.5
$ws.Range("I126").Value = 3998
$ws.Range("K126").Value = 11994
$ws.Range("M126").Value = -9524

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 984.8889
$ws.Range("J5").Value = 934.0909
$ws.Range("L5").Value = 2802.2727
$ws.Range("N5").Value = -3026.2727
# Row 34
$ws.Range("H34").Value = 66465.88
$ws.Range("J34").Value = 78399.28999999999
$ws.Range("L34").Value = 235197.87
$ws.Range("N34").Value = -235365.87
# Row 39
$ws.Range("H39").Value = 7409
$ws.Range("J39").Value = 7409
$ws.Range("L39").Value = 22227
$ws.Range("N39").Value = -22815
# Row 55
$ws.Range("H55").Value = 19998.5
$ws.Range("J55").Value = 19998.5
$ws.Range("L55").Value = 59995.5
$ws.Range("N55").Value = -60349.5
# Row 86
$ws.Range("H86").Value = 1185.7778
$ws.Range("I86").Value = 489.30768
$ws.Range("J86").Value = 2996.6
$ws.Range("K86").Value = 1467.92304
$ws.Range("L86").Value = 8989.799999999999
$ws.Range("M86").Value = -281.9230400000001
$ws.Range("N86").Value = -11361.8
# Row 89
$ws.Range("H89").Value = 1185.7778
$ws.Range("I89").Value = 489.30768
$ws.Range("J89").Value = 2996.6
$ws.Range("K89").Value = 4403.76912
$ws.Range("L89").Value = 26969.4
$ws.Range("M89").Value = 1524.23088
$ws.Range("N89").Value = -38825.39999999999
# Row 97
$ws.Range("H97").Value = 2519.7
$ws.Range("J97").Value = 2749.875
$ws.Range("L97").Value = 8249.625
$ws.Range("N97").Value = -9241.625
# Row 132
$ws.Range("H132").Value = 1267.091
$ws.Range("J132").Value = 1314
$ws.Range("L132").Value = 11826
$ws.Range("N132").Value = -16886
# Row 135
$ws.Range("H135").Value = 984.8889
$ws.Range("J135").Value = 934.0909
$ws.Range("L135").Value = 8406.8181
$ws.Range("N135").Value = -13476.8181
# Row 140
$ws.Range("H140").Value = 528921.1
$ws.Range("I140").Value = 528921.1
$ws.Range("K140").Value = 1586763.3
$ws.Range("M140").Value = -1581583.3

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 46
$ws.Range("H46").Value = 20681.941
$ws.Range("I46").Value = 4201.25
$ws.Range("K46").Value = 4201.25
$ws.Range("M46").Value = -4045.25
# Row 70
$ws.Range("H70").Value = 10000
$ws.Range("I70").Value = 10000
$ws.Range("K70").Value = 10000
$ws.Range("M70").Value = -9730
# Row 73
$ws.Range("H73").Value = 10000
$ws.Range("I73").Value = 10000
$ws.Range("K73").Value = 10000
$ws.Range("M73").Value = -9064
# Row 97
$ws.Range("H97").Value = 615.9167
$ws.Range("I97").Value = 644.63635
$ws.Range("J97").Value = 300
$ws.Range("K97").Value = 644.63635
$ws.Range("L97").Value = 300
$ws.Range("M97").Value = -148.63635
$ws.Range("N97").Value = -1292

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 511.4
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
# Row 22
$ws.Range("H22").Value = 2955.5715
$ws.Range("I22").Value = 3114.8333
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 3114.8333
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -2819.8333
$ws.Range("N22").Value = -2590
# Row 27
$ws.Range("H27").Value = 2955.5715
$ws.Range("I27").Value = 3114.8333
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 3114.8333
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -3007.8333
$ws.Range("N27").Value = -2214
# Row 46
$ws.Range("H46").Value = 2026.1428
$ws.Range("I46").Value = 1157.5
$ws.Range("J46").Value = 3184.3333
$ws.Range("K46").Value = 1157.5
$ws.Range("L46").Value = 3184.3333
$ws.Range("M46").Value = -969.5
$ws.Range("N46").Value = -3560.3333
# Row 68
$ws.Range("H68").Value = 4295.6665
$ws.Range("I68").Value = 4295.6665
$ws.Range("K68").Value = 4295.6665
$ws.Range("M68").Value = -3546.6665
# Row 71
$ws.Range("H71").Value = 4295.6665
$ws.Range("I71").Value = 4295.6665
$ws.Range("K71").Value = 21478.3325
$ws.Range("M71").Value = -17734.3325
# Row 122
$ws.Range("H122").Value = 2672.9412
$ws.Range("I122").Value = 2457.7693
$ws.Range("K122").Value = 7373.3079
$ws.Range("M122").Value = -4923.3079
# Row 132
$ws.Range("H132").Value = 2659.92
$ws.Range("I132").Value = 2793.75
$ws.Range("J132").Value = 2422
$ws.Range("K132").Value = 8381.25
$ws.Range("L132").Value = 7266
$ws.Range("M132").Value = -5851.25
$ws.Range("N132").Value = -12326

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 4
$ws.Range("H4").Value = 4999.5
$ws.Range("J4").Value = 4999.5
$ws.Range("L4").Value = 4999.5
$ws.Range("N4").Value = -5225.5
# Row 29
$ws.Range("H29").Value = 70145
$ws.Range("I29").Value = 70145
$ws.Range("K29").Value = 70145
$ws.Range("M29").Value = -69855
# Row 62
$ws.Range("H62").Value = 5483
$ws.Range("I62").Value = 5434.2
$ws.Range("J62").Value = 5564.3335
$ws.Range("K62").Value = 5434.2
$ws.Range("L62").Value = 5564.3335
$ws.Range("M62").Value = -4810.2
$ws.Range("N62").Value = -6812.3335
# Row 65
$ws.Range("H65").Value = 5483
$ws.Range("I65").Value = 5434.2
$ws.Range("J65").Value = 5564.3335
$ws.Range("K65").Value = 27171
$ws.Range("L65").Value = 27821.6675
$ws.Range("M65").Value = -24051
$ws.Range("N65").Value = -34061.6675
# Row 107
$ws.Range("H107").Value = 585.3570999999999
$ws.Range("I107").Value = 575.75
$ws.Range("J107").Value = 598.1667
$ws.Range("K107").Value = 1727.25
$ws.Range("L107").Value = 1794.5001
$ws.Range("M107").Value = 192.75
$ws.Range("N107").Value = -5634.5001
# Row 126
$ws.Range("H126").Value = 3149.25
$ws.Range("I126").Value = 1865.6666
$ws.Range("K126").Value = 5596.9998
$ws.Range("M126").Value = -3126.9998
# Row 132
$ws.Range("H132").Value = 9891.200000000001
$ws.Range("I132").Value = 3500.6667
$ws.Range("J132").Value = 19477
$ws.Range("K132").Value = 10502.0001
$ws.Range("L132").Value = 58431
$ws.Range("M132").Value = -7972.000100000001
$ws.Range("N132").Value = -63491

